# Updated Delivery Changes via GUI
# Row 17 values were edited in the sheet (Start Date, Store, Detail, Warehouse).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds a date-like string that must stay as plain text (it was stored
# as text before the edit too), so prefix with an apostrophe like a user typing
# into Excel would, to avoid Excel auto-converting it into a real date value.
$ws.Range("A17").Value = "'2025-04-26"
$ws.Range("B17").Value = "asd"
$ws.Range("C17").Value = "sdf"
$ws.Range("D17").Value = "dsfxgf"
